# Strategic Logistics Board - Warehouse Complaints
# Translate the English complaint "category" values (column D) to German,
# turn off the AutoFilter, and reset the sheet's scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map the English category labels to their German translations.
$categoryMap = @{
    "Damaged package" = "Verpackung beschädigt"
    "Too late"         = "Verspätet"
    "Missing parts"    = "Fehlende Teile"
    "Other"            = "Sonstiges"
}

$lastRow = 200
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $current = $cell.Value()
    if ($categoryMap.ContainsKey($current)) {
        $cell.Value = $categoryMap[$current]
    }
}

# Remove the AutoFilter from the sheet.
if ($ws.AutoFilterMode()) {
    $ws.AutoFilterMode = $false
}

# Reset the view: scroll back to the top and move the selection to K5.
$ws.Range("A1").Select() | Out-Null
$ws.Range("K5").Select() | Out-Null
